$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the extra trailing row (old row 16: B64/D=PCLOUD/E=2 duplicate) ---
$ws.Rows(16).Delete()

# --- Remove column E (constant "1" column); old F/G shift left to E/F ---
$ws.Columns("E").Delete()

# --- New column F (previously G) gets a 3-decimal numeric format ---
$ws.Range("F1:F15").NumberFormat = "0.000"
$ws.Columns("F").ColumnWidth = 10

# --- Update the view selection to cover the whole used range ---
$ws.Range("A1:F15").Select()

# --- Conditional formatting on column D: highlight rows containing each label ---
$rng = $ws.Range("D1:D1048576")

$fcICOLOR = $rng.FormatConditions.Add(9)
$fcICOLOR.TextOperator = 0
$fcICOLOR.Text = "ICOLOR"
$fcICOLOR.Formula1 = 'NOT(ISERROR(SEARCH("ICOLOR",D1)))'
$fcICOLOR.Font.Color = 393372
$fcICOLOR.Interior.Color = 13551615

$fcIDEPTH = $rng.FormatConditions.Add(9)
$fcIDEPTH.TextOperator = 0
$fcIDEPTH.Text = "IDEPTH"
$fcIDEPTH.Formula1 = 'NOT(ISERROR(SEARCH("IDEPTH",D1)))'
$fcIDEPTH.Font.Color = 22428
$fcIDEPTH.Interior.Color = 10284031
$fcIDEPTH.SetFirstPriority()

$fcPCLOUD = $rng.FormatConditions.Add(9)
$fcPCLOUD.TextOperator = 0
$fcPCLOUD.Text = "PCLOUD"
$fcPCLOUD.Formula1 = 'NOT(ISERROR(SEARCH("PCLOUD",D1)))'
$fcPCLOUD.Font.Color = 24832
$fcPCLOUD.Interior.Color = 13561798
$fcPCLOUD.SetFirstPriority()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1
